$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the "Directory" path values (column D) for the 4 data rows:
# old: data/Multimedia_Data/Documentation/  ->  new: data/multimedia/documentation/
$ws.Range("D2:D5").Value = "data/multimedia/documentation/"

# Update the active selection to match the authored state
$ws.Range("D8").Select()
